$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values keep their exact text representation
# (avoids Excel auto-converting numeric-looking strings and dropping
# trailing zeros / collapsing dotted values).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.174.09"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.915.05"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "0.7390"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").Value = "244.15"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.3130"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").Value = "26.82"
$ws.Range("D10").Value = "0.06974"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Value = "0.7811"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "0.07991"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.904.34"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "5.291"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "92.46"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "30.181.05"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "5.928"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").Value = "242.75"
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "0.000007839"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.147.07"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "7.163"
$ws.Range("E24").Value = "  +7.40%  "
$ws.Range("D25").Value = "9.439"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "168.84"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").Value = "19.13"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "0.1285"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").Value = "2.075"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").Value = "1.353"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "4.348"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").Value = "4.113"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "0.05174"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "1.303"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "2.803"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "75.22"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "0.4519"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "1.966"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "7.895"
$ws.Range("E44").Value = "  +5.60%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "0.8396"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "9.946"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "101.63"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").Value = "37.30"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "2.053.84"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "0.1199"
$ws.Range("E51").Value = "  +2.06%  "
